# Update "Förändrad" (changed) date column for rows 2-6 to the new date serial 45204 (2023-10-05)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C6").Value = 45204
